$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value2 = $text
    $r.ClearFormats()
}

# Qty executed upto date (plain numeric cells)
$ws.Range("C8").Value = 7
$ws.Range("C9").Value = 60
$ws.Range("C10").Value = 11
$ws.Range("C11").Value = 49
$ws.Range("C12").Value = 45
$ws.Range("C13").Value = 83
$ws.Range("C14").Value = 88
$ws.Range("C15").Value = 71
$ws.Range("C16").Value = 55
$ws.Range("C17").Value = 37

# Upto date Amount cells (stored as text strings with a fixed "x.00" format)
Set-TextValue "G9" "15360.00"
Set-TextValue "G10" "5192.00"
Set-TextValue "G11" "32438.00"
Set-TextValue "G13" "11288.00"
Set-TextValue "G14" "2024.00"

# Grand totals
Set-TextValue "G19" "66302.00"
Set-TextValue "H19" "66302.00"
Set-TextValue "G21" "66302.00"
Set-TextValue "H21" "66302.00"
